$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header cell in H1 -- copy formatting from the neighboring
# header cell (G1) so it matches the existing bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for data rows 2-15 (plain numeric, unstyled like
# the other numeric columns).
$saveValues = @(1, 0, 1, 1, 1, 1, 0, 1, 0, 1, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
